# Horarios Línea 141 - actualización 558
# Actualiza la hora de scrapeo (04:52:35 -> 05:16:08) y recalcula las
# filas de arribos para las tres hojas del libro.

$wb = $excel.ActiveWorkbook

$oldScrap = "04:52:35"
$newScrap = "05:16:08"

# ---------------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newScrap"
$ws1.Range("A3").Value = "Total filas: 13"

$sheet1Rows = @(
    @("05:17", "17_ROMERO",      6),
    @("05:22", "23_HERNANDEZ",   6),
    @("05:44", "14_ABASTO",      28),
    @("05:47", "17_ROMERO",      31),
    @("06:01", "16_SANTA ANA",   45),
    @("06:09", "10_OLMOS",       53),
    @("06:16", "215A_EL PATO",   60),
    @("06:30", "23_HERNANDEZ",   74),
    @("06:34", "11_ETCHEVERRY",  78),
    @("06:39", "17X38_ROMERO",   83),
    @("06:41", "16_SANTA ANA",   85),
    @("06:57", "215A_EL PATO",   101),
    @("06:59", "225_GOMEZ",      103)
)

$r = 6
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $newScrap
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = "LP1912"
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newScrap"
$ws2.Range("A3").Value = "Total filas: 2"

$sheet2Rows = @(
    @("06:16", "215A_EL PATO", 60),
    @("06:57", "215A_EL PATO", 101)
)

$r = 6
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $newScrap
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = "LP1912"
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Hoja 3: 6203-6173 (sin arribos, solo se actualiza el timestamp)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newScrap"
